$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is referenced from the Overview sheet (B/C columns)
#    as well as the "Status" column (C) on each per-locale sheet.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Per-locale sheets ("zh-cn", "de-de"): populate the "Latest Target File"
#    (F) and "Latest Handback File" (G) columns for rows 2 & 3, update the
#    Status column (C) and stamp the "Latest Handback DateTime" column (H)
#    with the handback timestamp.
# ---------------------------------------------------------------------------
$localeInfo = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-03-20 16:15:35" },
    @{ Name = "de-de"; HandbackTime = "2016-03-20 16:15:40" }
)

foreach ($info in $localeInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Snapshot the existing hyperlink targets (and their display text) keyed
    # by the cell address they live on, so the new F/G hyperlinks can reuse
    # the exact same target URLs as the existing A/D hyperlinks.
    $linkAddress = @{}
    $linkDisplay = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        $linkAddress[$addr] = $hl.Address
        $linkDisplay[$addr] = $hl.TextToDisplay
    }

    foreach ($row in @(2, 3)) {
        $aAddr = "`$A`$" + $row
        $dAddr = "`$D`$" + $row

        $fCell = $ws.Range("F" + $row)
        $gCell = $ws.Range("G" + $row)

        $fCell.Value = $linkDisplay[$aAddr]
        $fCell.Style = "HyperLink"
        $ws.Hyperlinks.Add($fCell, $linkAddress[$aAddr], "", "", $linkDisplay[$aAddr]) | Out-Null

        $gCell.Value = $linkDisplay[$dAddr]
        $gCell.Style = "HyperLink"
        $ws.Hyperlinks.Add($gCell, $linkAddress[$dAddr], "", "", $linkDisplay[$dAddr]) | Out-Null

        # Status column
        $ws.Range("C" + $row).Value = $newStatus

        # Latest Handback DateTime column
        $ws.Range("H" + $row).Value = $info.HandbackTime
    }
}

Write-Host "Report for handback generated"
